$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new test case row above the current row 2 (TC_01),
#    pushing the existing TC_01/TC_02/TC_03 rows down to rows 3-5.
# ------------------------------------------------------------------
$ws.Rows.Item(2).Insert()

# Copy the formatting of the (now) row 3 into the new row 2 so it
# starts out with the same border / wrap / vertical-top look as the
# rest of the table.
$ws.Range("A3:E3").Copy()
$ws.Range("A2:E2").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2. Fill in the new test case (row 2).
# ------------------------------------------------------------------
$ws.Cells.Item(2, 1).Value = "TC_01"
$ws.Cells.Item(2, 3).Value = "1) Fetch the page title`n2) Compare it with the expected page title"
$ws.Cells.Item(2, 4).Value = '1) page title should be "IMDB Top 250 - Imdb"'
$ws.Cells.Item(2, 2).Value = "Verift title and the page"
$ws.Cells.Item(2, 5).Value = "Pass"

# The Status cell for the new row does not wrap like the rest of the row.
$ws.Range("E2").WrapText = $false

$ws.Rows.Item(2).RowHeight = 45

# ------------------------------------------------------------------
# 3. Renumber the Test Case IDs that got shifted down.
# ------------------------------------------------------------------
$ws.Cells.Item(3, 1).Value = "TC_02"
$ws.Cells.Item(4, 1).Value = "TC_03"
$ws.Cells.Item(5, 1).Value = "TC_04"

# ------------------------------------------------------------------
# 4. Formatting tweaks.
# ------------------------------------------------------------------
# Column A (Test Case ID) - align top & wrap for the data rows.
$ws.Range("A4:A5").VerticalAlignment = -4160
$ws.Range("A4:A5").WrapText = $true

# Column D (Expected result) - align top, keep the header from wrapping.
$ws.Range("D1:D5").VerticalAlignment = -4160
$ws.Range("D1").WrapText = $false
$ws.Range("D3:D5").WrapText = $true

# Column D is now wider to fit the new content.
$ws.Columns.Item(4).ColumnWidth = 24

# ------------------------------------------------------------------
# 5. Reset the selection back to A1.
# ------------------------------------------------------------------
$ws.Range("A1").Select()
